{"js": "const replacements = [\n  [\"16\u00f76=\", \"23\u00f79=\"],\n  [\"57\u00f79=\", \"71\u00f79=\"],\n  [\"64\u00f77=\", \"73\u00f76=\"],\n  [\"79\u00f76=\", \"72\u00f73=\"],\n  [\"29\u00f78=\", \"40\u00f72=\"],\n  [\"80\u00f76=\", \"90\u00f74=\"],\n  [\"33\u00f74=\", \"68\u00f73=\"],\n  [\"57\u00f72=\", \"39\u00f79=\"],\n  [\"68\u00f78=\", \"46\u00f72=\"],\n  [\"32\u00f77=\", \"52\u00f77=\"],\n  [\"65\u00f76=\", \"57\u00f76=\"],\n  [\"14\u00f73=\", \"97\u00f74=\"],\n  [\"54\u00f79=\", \"43\u00f79=\"],\n  [\"74\u00f76=\", \"84\u00f74=\"],\n  [\"61\u00f74=\", \"10\u00f75=\"],\n  [\"41\u00f78=\", \"84\u00f74=\"],\n  [\"79\u00f74=\", \"60\u00f79=\"],\n  [\"74\u00f77=\", \"91\u00f72=\"],\n  [\"40\u00f75=\", \"44\u00f72=\"],\n  [\"77\u00f79=\", \"10\u00f73=\"],\n  [\"86\u00f76=\", \"99\u00f75=\"],\n  [\"37\u00f77=\", \"50\u00f75=\"],\n  [\"83\u00f75=\", \"98\u00f74=\"],\n  [\"39\u00f75=\", \"90\u00f74=\"],\n  [\"83\u00f79=\", \"93\u00f74=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load('items');\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"16\u00f76=\", \"23\u00f79=\"),\n    @(\"57\u00f79=\", \"71\u00f79=\"),\n    @(\"64\u00f77=\", \"73\u00f76=\"),\n    @(\"79\u00f76=\", \"72\u00f73=\"),\n    @(\"29\u00f78=\", \"40\u00f72=\"),\n    @(\"80\u00f76=\", \"90\u00f74=\"),\n    @(\"33\u00f74=\", \"68\u00f73=\"),\n    @(\"57\u00f72=\", \"39\u00f79=\"),\n    @(\"68\u00f78=\", \"46\u00f72=\"),\n    @(\"32\u00f77=\", \"52\u00f77=\"),\n    @(\"65\u00f76=\", \"57\u00f76=\"),\n    @(\"14\u00f73=\", \"97\u00f74=\"),\n    @(\"54\u00f79=\", \"43\u00f79=\"),\n    @(\"74\u00f76=\", \"84\u00f74=\"),\n    @(\"61\u00f74=\", \"10\u00f75=\"),\n    @(\"41\u00f78=\", \"84\u00f74=\"),\n    @(\"79\u00f74=\", \"60\u00f79=\"),\n    @(\"74\u00f77=\", \"91\u00f72=\"),\n    @(\"40\u00f75=\", \"44\u00f72=\"),\n    @(\"77\u00f79=\", \"10\u00f73=\"),\n    @(\"86\u00f76=\", \"99\u00f75=\"),\n    @(\"37\u00f77=\", \"50\u00f75=\"),\n    @(\"83\u00f75=\", \"98\u00f74=\"),\n    @(\"39\u00f75=\", \"90\u00f74=\"),\n    @(\"83\u00f79=\", \"93\u00f74=\"),\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $newText\n    $find.Execute($oldText, $false, $true, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null\n}"}
